$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H7").Value = 12250
$ws_ALC.Range("I7").Value = 0
$ws_ALC.Range("J7").Value = 12250
$ws_ALC.Range("K7").Value = 0
$ws_ALC.Range("L7").Value = 12250
$ws_ALC.Range("M7").ClearContents()
$ws_ALC.Range("N7").Value = -12474

$ws_ALC.Range("H14").Value = 12250
$ws_ALC.Range("I14").Value = 0
$ws_ALC.Range("J14").Value = 12250
$ws_ALC.Range("K14").Value = 0
$ws_ALC.Range("L14").Value = 12250
$ws_ALC.Range("M14").ClearContents()
$ws_ALC.Range("N14").Value = -12632

$ws_ALC.Range("H74").Value = 5436.4546
$ws_ALC.Range("I74").Value = 4033.6667
$ws_ALC.Range("K74").Value = 4033.6667
$ws_ALC.Range("M74").Value = -3097.6667

$ws_ALC.Range("H77").Value = 5436.4546
$ws_ALC.Range("I77").Value = 4033.6667
$ws_ALC.Range("K77").Value = 20168.3335
$ws_ALC.Range("M77").Value = -15488.3335

$ws_ALC.Range("H86").Value = 4050.6155
$ws_ALC.Range("I86").Value = 2429.5334
$ws_ALC.Range("J86").Value = 6261.1816
$ws_ALC.Range("K86").Value = 2429.5334
$ws_ALC.Range("L86").Value = 6261.1816
$ws_ALC.Range("M86").Value = -1306.5334
$ws_ALC.Range("N86").Value = -8507.1816

$ws_ALC.Range("H89").Value = 4050.6155
$ws_ALC.Range("I89").Value = 2429.5334
$ws_ALC.Range("J89").Value = 6261.1816
$ws_ALC.Range("K89").Value = 12147.667
$ws_ALC.Range("L89").Value = 31305.908
$ws_ALC.Range("M89").Value = -6531.666999999999
$ws_ALC.Range("N89").Value = -42537.908

$ws_ALC.Range("H129").Value = 2809.3137
$ws_ALC.Range("J129").Value = 909.381
$ws_ALC.Range("L129").Value = 2728.143
$ws_ALC.Range("N129").Value = -12728.143

$ws_ALC.Range("H131").Value = 4202.4165
$ws_ALC.Range("I131").Value = 2651.25
$ws_ALC.Range("J131").Value = 4978
$ws_ALC.Range("K131").Value = 7953.75
$ws_ALC.Range("L131").Value = 14934
$ws_ALC.Range("M131").Value = -2913.75
$ws_ALC.Range("N131").Value = -25014

$ws_ALC.Range("H135").Value = 1836.6
$ws_ALC.Range("I135").Value = 1836.6
$ws_ALC.Range("J135").Value = 0
$ws_ALC.Range("K135").Value = 16529.4
$ws_ALC.Range("L135").Value = 0
$ws_ALC.Range("M135").Value = -13994.4
$ws_ALC.Range("N135").ClearContents()

$ws_ALC.Range("H137").Value = 1392.0588
$ws_ALC.Range("I137").Value = 1134.5454
$ws_ALC.Range("J137").Value = 1587.4138
$ws_ALC.Range("K137").Value = 3403.6362
$ws_ALC.Range("L137").Value = 4762.2414
$ws_ALC.Range("M137").Value = -853.6361999999999
$ws_ALC.Range("N137").Value = -9862.241399999999

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 2358.8064
$ws_ARM.Range("I61").Value = 1676.5
$ws_ARM.Range("K61").Value = 1676.5
$ws_ARM.Range("M61").Value = -1464.5

$ws_ARM.Range("H74").Value = 2591.7036
$ws_ARM.Range("I74").Value = 1947.2307
$ws_ARM.Range("K74").Value = 1947.2307
$ws_ARM.Range("M74").Value = -1073.2307

$ws_ARM.Range("H77").Value = 2591.7036
$ws_ARM.Range("I77").Value = 1947.2307
$ws_ARM.Range("K77").Value = 9736.1535
$ws_ARM.Range("M77").Value = -5368.1535

$ws_ARM.Range("H110").Value = 25051874
$ws_ARM.Range("I110").Value = 31313698
$ws_ARM.Range("K110").Value = 31313698
$ws_ARM.Range("M110").Value = -31311653

$ws_ARM.Range("H136").Value = 2358.8064
$ws_ARM.Range("I136").Value = 1676.5
$ws_ARM.Range("K136").Value = 5029.5
$ws_ARM.Range("M136").Value = -2479.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H12").Value = 1573.75
$ws_BSM.Range("I12").Value = 1573.75
$ws_BSM.Range("K12").Value = 1573.75
$ws_BSM.Range("M12").Value = -1405.75

$ws_BSM.Range("H94").Value = 200163.4
$ws_BSM.Range("I94").Value = 200163.4
$ws_BSM.Range("K94").Value = 200163.4
$ws_BSM.Range("M94").Value = -199712.4

$ws_BSM.Range("H99").Value = 1682.8572
$ws_BSM.Range("I99").Value = 1624.2858
$ws_BSM.Range("J99").Value = 1741.4286
$ws_BSM.Range("K99").Value = 1624.2858
$ws_BSM.Range("L99").Value = 1741.4286
$ws_BSM.Range("M99").Value = -126.2858000000001
$ws_BSM.Range("N99").Value = -4737.4286

$ws_BSM.Range("H134").Value = 2396.0312
$ws_BSM.Range("I134").Value = 2443.8965
$ws_BSM.Range("J134").Value = 1933.3334
$ws_BSM.Range("K134").Value = 7331.689499999999
$ws_BSM.Range("L134").Value = 5800.0002
$ws_BSM.Range("M134").Value = -4796.689499999999
$ws_BSM.Range("N134").Value = -10870.0002

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H99").Value = 7740.619
$ws_CRP.Range("I99").Value = 2945.5833
$ws_CRP.Range("K99").Value = 2945.5833
$ws_CRP.Range("M99").Value = -1447.5833

$ws_CRP.Range("H107").Value = 4210.2856
$ws_CRP.Range("I107").Value = 6936.3125
$ws_CRP.Range("J107").Value = 575.5833
$ws_CRP.Range("K107").Value = 6936.3125
$ws_CRP.Range("L107").Value = 575.5833
$ws_CRP.Range("M107").Value = -5016.3125
$ws_CRP.Range("N107").Value = -4415.5833

$ws_CRP.Range("H126").Value = 7740.619
$ws_CRP.Range("I126").Value = 2945.5833
$ws_CRP.Range("K126").Value = 8836.749899999999
$ws_CRP.Range("M126").Value = -6366.749899999999

$ws_CRP.Range("H134").Value = 1749.6666
$ws_CRP.Range("I134").Value = 1248.6154
$ws_CRP.Range("K134").Value = 3745.8462
$ws_CRP.Range("M134").Value = -1210.8462

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H34").Value = 1924.875
$ws_CUL.Range("J34").Value = 2516.5
$ws_CUL.Range("L34").Value = 7549.5
$ws_CUL.Range("N34").Value = -7717.5

$ws_CUL.Range("H54").Value = 1000
$ws_CUL.Range("J54").Value = 1000
$ws_CUL.Range("L54").Value = 3000
$ws_CUL.Range("N54").Value = -4118

$ws_CUL.Range("H68").Value = 1763.9375
$ws_CUL.Range("I68").Value = 848.4074000000001
$ws_CUL.Range("J68").Value = 2432.027
$ws_CUL.Range("K68").Value = 2545.2222
$ws_CUL.Range("L68").Value = 7296.081
$ws_CUL.Range("M68").Value = -1734.2222
$ws_CUL.Range("N68").Value = -8918.081

$ws_CUL.Range("H71").Value = 1763.9375
$ws_CUL.Range("I71").Value = 848.4074000000001
$ws_CUL.Range("J71").Value = 2432.027
$ws_CUL.Range("K71").Value = 7635.6666
$ws_CUL.Range("L71").Value = 21888.243
$ws_CUL.Range("M71").Value = -3579.6666
$ws_CUL.Range("N71").Value = -30000.243

$ws_CUL.Range("H140").Value = 4882.3447
$ws_CUL.Range("I140").Value = 5976.9
$ws_CUL.Range("J140").Value = 2450
$ws_CUL.Range("K140").Value = 17930.7
$ws_CUL.Range("L140").Value = 7350
$ws_CUL.Range("M140").Value = -12750.7
$ws_CUL.Range("N140").Value = -17710

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H12").Value = 4689312.5
$ws_GSM.Range("I12").Value = 5001600
$ws_GSM.Range("J12").Value = 5000
$ws_GSM.Range("K12").Value = 5001600
$ws_GSM.Range("L12").Value = 5000
$ws_GSM.Range("M12").Value = -5001460
$ws_GSM.Range("N12").Value = -5280

$ws_GSM.Range("H21").Value = 5000
$ws_GSM.Range("I21").Value = 5000
$ws_GSM.Range("J21").Value = 0
$ws_GSM.Range("K21").Value = 5000
$ws_GSM.Range("L21").Value = 0
$ws_GSM.Range("M21").Value = -4827
$ws_GSM.Range("N21").ClearContents()

$ws_GSM.Range("H30").Value = 5000
$ws_GSM.Range("I30").Value = 5000
$ws_GSM.Range("J30").Value = 0
$ws_GSM.Range("K30").Value = 5000
$ws_GSM.Range("L30").Value = 0
$ws_GSM.Range("M30").Value = -4895
$ws_GSM.Range("N30").ClearContents()

$ws_GSM.Range("H38").Value = 0
$ws_GSM.Range("J38").Value = 0
$ws_GSM.Range("L38").Value = 0
$ws_GSM.Range("N38").ClearContents()

$ws_GSM.Range("H40").Value = 8000
$ws_GSM.Range("J40").Value = 8000
$ws_GSM.Range("L40").Value = 8000
$ws_GSM.Range("N40").Value = -8302

$ws_GSM.Range("H102").Value = 263280.2
$ws_GSM.Range("I102").Value = 2876.7334
$ws_GSM.Range("J102").Value = 751536.6
$ws_GSM.Range("K102").Value = 2876.7334
$ws_GSM.Range("L102").Value = 751536.6
$ws_GSM.Range("M102").Value = -1254.7334
$ws_GSM.Range("N102").Value = -754780.6

$ws_GSM.Range("H122").Value = 2713.818
$ws_GSM.Range("I122").Value = 2474.125
$ws_GSM.Range("J122").Value = 3353
$ws_GSM.Range("K122").Value = 7422.375
$ws_GSM.Range("L122").Value = 10059
$ws_GSM.Range("M122").Value = -4972.375
$ws_GSM.Range("N122").Value = -14959

$ws_GSM.Range("H132").Value = 4287.4443
$ws_GSM.Range("I132").Value = 4369.7144
$ws_GSM.Range("J132").Value = 3999.5
$ws_GSM.Range("K132").Value = 13109.1432
$ws_GSM.Range("L132").Value = 11998.5
$ws_GSM.Range("M132").Value = -10579.1432
$ws_GSM.Range("N132").Value = -17058.5

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 1497.8
$ws_LTW.Range("I7").Value = 1164.1111
$ws_LTW.Range("J7").Value = 1998.3334
$ws_LTW.Range("K7").Value = 1164.1111
$ws_LTW.Range("L7").Value = 1998.3334
$ws_LTW.Range("M7").Value = -1052.1111
$ws_LTW.Range("N7").Value = -2222.3334

$ws_LTW.Range("H19").Value = 9929.143
$ws_LTW.Range("I19").Value = 3000
$ws_LTW.Range("J19").Value = 12700.8
$ws_LTW.Range("K19").Value = 3000
$ws_LTW.Range("L19").Value = 12700.8
$ws_LTW.Range("M19").Value = -2830
$ws_LTW.Range("N19").Value = -13040.8

$ws_LTW.Range("H126").Value = 1497.8
$ws_LTW.Range("I126").Value = 1164.1111
$ws_LTW.Range("J126").Value = 1998.3334
$ws_LTW.Range("K126").Value = 3492.3333
$ws_LTW.Range("L126").Value = 5995.0002
$ws_LTW.Range("M126").Value = -1022.3333
$ws_LTW.Range("N126").Value = -10935.0002

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H122").Value = 1002.6667
$ws_WVR.Range("I122").Value = 1004.8571
$ws_WVR.Range("J122").Value = 995
$ws_WVR.Range("K122").Value = 3014.5713
$ws_WVR.Range("L122").Value = 2985
$ws_WVR.Range("M122").Value = -564.5712999999996
$ws_WVR.Range("N122").Value = -7885
